# feat: read data from excel to list
#
# The "Gift" sheet tracks quantity of gifts per id; this adds the missing
# "remaining" value (column D) for the 4th row (id 4, "Gấu bông") and makes
# the "Gift" sheet the active/selected sheet & cell, as left by the author
# after entering the value.

$wb = $excel.ActiveWorkbook

$giftSheet = $wb.Worksheets.Item("Gift")

# Fill in the missing "remaining" value for row 5 (id 4)
$giftSheet.Range("D5").Value = 4

# Author ended up with the Gift sheet active and D6 selected
$giftSheet.Activate()
$giftSheet.Range("D6").Select() | Out-Null
